{"js": "// Update the two-digit multiplication problems in the document's table.\n// Each old expression is unique in the document, so a search-and-replace\n// keyed on the exact old text safely targets the correct run while\n// preserving all existing run formatting (fonts, size, etc.).\nconst replacements = [\n  [\"14\u00d725=\", \"91\u00d772=\"],\n  [\"93\u00d731=\", \"14\u00d747=\"],\n  [\"88\u00d776=\", \"61\u00d764=\"],\n  [\"70\u00d733=\", \"51\u00d728=\"],\n  [\"69\u00d754=\", \"22\u00d798=\"],\n  [\"36\u00d714=\", \"16\u00d728=\"],\n  [\"13\u00d779=\", \"51\u00d722=\"],\n  [\"79\u00d729=\", \"18\u00d747=\"],\n  [\"14\u00d758=\", \"61\u00d795=\"],\n  [\"30\u00d757=\", \"21\u00d731=\"],\n  [\"62\u00d774=\", \"65\u00d727=\"],\n  [\"92\u00d732=\", \"94\u00d770=\"],\n  [\"90\u00d781=\", \"18\u00d739=\"],\n  [\"68\u00d723=\", \"31\u00d737=\"],\n  [\"31\u00d720=\", \"12\u00d780=\"],\n  [\"98\u00d763=\", \"58\u00d799=\"],\n  [\"25\u00d776=\", \"95\u00d731=\"],\n  [\"68\u00d783=\", \"99\u00d731=\"],\n  [\"11\u00d731=\", \"57\u00d735=\"],\n  [\"56\u00d767=\", \"80\u00d752=\"],\n  [\"36\u00d797=\", \"97\u00d743=\"],\n  [\"74\u00d773=\", \"72\u00d759=\"],\n  [\"82\u00d725=\", \"28\u00d734=\"],\n  [\"17\u00d790=\", \"30\u00d755=\"],\n  [\"74\u00d775=\", \"99\u00d747=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"14\u00d725=\", \"91\u00d772=\"),\n    @(\"93\u00d731=\", \"14\u00d747=\"),\n    @(\"88\u00d776=\", \"61\u00d764=\"),\n    @(\"70\u00d733=\", \"51\u00d728=\"),\n    @(\"69\u00d754=\", \"22\u00d798=\"),\n    @(\"36\u00d714=\", \"16\u00d728=\"),\n    @(\"13\u00d779=\", \"51\u00d722=\"),\n    @(\"79\u00d729=\", \"18\u00d747=\"),\n    @(\"14\u00d758=\", \"61\u00d795=\"),\n    @(\"30\u00d757=\", \"21\u00d731=\"),\n    @(\"62\u00d774=\", \"65\u00d727=\"),\n    @(\"92\u00d732=\", \"94\u00d770=\"),\n    @(\"90\u00d781=\", \"18\u00d739=\"),\n    @(\"68\u00d723=\", \"31\u00d737=\"),\n    @(\"31\u00d720=\", \"12\u00d780=\"),\n    @(\"98\u00d763=\", \"58\u00d799=\"),\n    @(\"25\u00d776=\", \"95\u00d731=\"),\n    @(\"68\u00d783=\", \"99\u00d731=\"),\n    @(\"11\u00d731=\", \"57\u00d735=\"),\n    @(\"56\u00d767=\", \"80\u00d752=\"),\n    @(\"36\u00d797=\", \"97\u00d743=\"),\n    @(\"74\u00d773=\", \"72\u00d759=\"),\n    @(\"82\u00d725=\", \"28\u00d734=\"),\n    @(\"17\u00d790=\", \"30\u00d755=\"),\n    @(\"74\u00d775=\", \"99\u00d747=\")\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $p[0]\n    $find.Replacement.Text = $p[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
